# Insert a new data row at row 52 (pushing the existing rows 52-181 down to 53-182),
# then populate the newly inserted row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 52, shifting rows 52-181 -> 53-182.
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with the new record's data.
$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "Vega Modelo de Temuco"
$ws.Range("C52").Value = "La Araucanía"
$ws.Range("D52").Value = 44525
$ws.Range("E52").Value = 9
$ws.Range("F52").Value = 100112039
$ws.Range("G52").Value = "Ciboulette"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 30
$ws.Range("K52").Value = 5000
$ws.Range("L52").Value = 5000
$ws.Range("M52").Value = 5000
$ws.Range("N52").Value = "`$/docena de atados"
$ws.Range("O52").Value = "Provincia de Cautín"
$ws.Range("P52").Value = 1667
$ws.Range("Q52").Value = 3
$ws.Range("R52").Value = "Hortaliza"

# Keep the date-formatted style consistent with the rest of column D.
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
